$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 29000
$ws.Range("I26").Value = 29000
$ws.Range("K26").Value = 29000
$ws.Range("M26").Value = -28656
$ws.Range("H40").Value = 3433.2727
$ws.Range("J40").Value = 3696
$ws.Range("L40").Value = 3696
$ws.Range("N40").Value = -4046
$ws.Range("H69").Value = 30002.334
$ws.Range("I69").Value = 9996
$ws.Range("J69").Value = 70015
$ws.Range("K69").Value = 29988
$ws.Range("L69").Value = 210045
$ws.Range("M69").Value = -29114
$ws.Range("N69").Value = -211793
$ws.Range("H72").Value = 30002.334
$ws.Range("I72").Value = 9996
$ws.Range("J72").Value = 70015
$ws.Range("K72").Value = 89964
$ws.Range("L72").Value = 630135
$ws.Range("M72").Value = -85596
$ws.Range("N72").Value = -638871
$ws.Range("H76").Value = 4119.8
$ws.Range("I76").Value = 3900.4285
$ws.Range("J76").Value = 4631.6665
$ws.Range("K76").Value = 3900.4285
$ws.Range("L76").Value = 4631.6665
$ws.Range("M76").Value = -3585.4285
$ws.Range("N76").Value = -5261.6665
$ws.Range("H79").Value = 4119.8
$ws.Range("I79").Value = 3900.4285
$ws.Range("J79").Value = 4631.6665
$ws.Range("K79").Value = 3900.4285
$ws.Range("L79").Value = 4631.6665
$ws.Range("M79").Value = -2808.4285
$ws.Range("N79").Value = -6815.6665
$ws.Range("H87").Value = 124940
$ws.Range("J87").Value = 124940
$ws.Range("L87").Value = 124940
$ws.Range("N87").Value = -127436
$ws.Range("H88").Value = 899155.7
$ws.Range("I88").Value = 1899.6666
$ws.Range("K88").Value = 1899.6666
$ws.Range("M88").Value = -1493.6666
$ws.Range("H90").Value = 124940
$ws.Range("J90").Value = 124940
$ws.Range("L90").Value = 374820
$ws.Range("N90").Value = -387300
$ws.Range("H91").Value = 899155.7
$ws.Range("I91").Value = 1899.6666
$ws.Range("K91").Value = 1899.6666
$ws.Range("M91").Value = -495.6666
$ws.Range("H103").Value = 674.6875
$ws.Range("I103").Value = 468.33334
$ws.Range("K103").Value = 1405.00002
$ws.Range("M103").Value = -819.0000199999999
$ws.Range("H132").Value = 2401.0527
$ws.Range("I132").Value = 2413.875
$ws.Range("K132").Value = 7241.625
$ws.Range("M132").Value = -4711.625
$ws.Range("H141").Value = 2451.8
$ws.Range("I141").Value = 2451.8
$ws.Range("K141").Value = 7355.400000000001
$ws.Range("M141").Value = -2175.400000000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 316.3158
$ws.Range("I5").Value = 235.66667
$ws.Range("K5").Value = 235.66667
$ws.Range("M5").Value = -123.66667
$ws.Range("H32").Value = 23816252
$ws.Range("I32").Value = 26321016
$ws.Range("K32").Value = 26321016
$ws.Range("M32").Value = -26320729
$ws.Range("H34").Value = 265247.5
$ws.Range("J34").Value = 265247.5
$ws.Range("L34").Value = 265247.5
$ws.Range("N34").Value = -265789.5
$ws.Range("H42").Value = 5000
$ws.Range("I42").Value = 5000
$ws.Range("K42").Value = 5000
$ws.Range("M42").Value = -4514
$ws.Range("H88").Value = 1314.5333
$ws.Range("I88").Value = 1445.25
$ws.Range("K88").Value = 1445.25
$ws.Range("M88").Value = -1039.25
$ws.Range("H91").Value = 1314.5333
$ws.Range("I91").Value = 1445.25
$ws.Range("K91").Value = 1445.25
$ws.Range("M91").Value = -41.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 316.3158
$ws.Range("I4").Value = 235.66667
$ws.Range("K4").Value = 235.66667
$ws.Range("M4").Value = -120.66667
$ws.Range("H94").Value = 858.95654
$ws.Range("I94").Value = 939.8421
$ws.Range("K94").Value = 939.8421
$ws.Range("M94").Value = -488.8421
$ws.Range("H105").Value = 2116.6365
$ws.Range("I105").Value = 1809.2222
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 1809.2222
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -62.22219999999993
$ws.Range("N105").Value = -6994
$ws.Range("H134").Value = 47698.652
$ws.Range("I134").Value = 868.5714
$ws.Range("J134").Value = 120545.445
$ws.Range("K134").Value = 2605.7142
$ws.Range("L134").Value = 361636.335
$ws.Range("M134").Value = -70.71420000000035
$ws.Range("N134").Value = -366706.335

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 12000
$ws.Range("J3").Value = 12000
$ws.Range("L3").Value = 12000
$ws.Range("N3").Value = -12226
$ws.Range("H58").Value = 1308.909
$ws.Range("I58").Value = 1350
$ws.Range("K58").Value = 1350
$ws.Range("M58").Value = -1147
$ws.Range("H117").Value = 87501
$ws.Range("J117").Value = 87501
$ws.Range("L117").Value = 87501
$ws.Range("N117").Value = -96679
$ws.Range("H136").Value = 1308.909
$ws.Range("I136").Value = 1350
$ws.Range("K136").Value = 4050
$ws.Range("M136").Value = -1500

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 40000
$ws.Range("I18").Value = 40000
$ws.Range("K18").Value = 40000
$ws.Range("M18").Value = -39707
$ws.Range("H113").Value = 2980.5
$ws.Range("I113").Value = 1011
$ws.Range("J113").Value = 4950
$ws.Range("K113").Value = 1011
$ws.Range("L113").Value = 4950
$ws.Range("M113").Value = 1159
$ws.Range("N113").Value = -9290
$ws.Range("H122").Value = 1725.9286
$ws.Range("I122").Value = 1698.091
$ws.Range("K122").Value = 5094.272999999999
$ws.Range("M122").Value = -2644.272999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3249.5
$ws.Range("J68").Value = 3500
$ws.Range("L68").Value = 3500
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 3249.5
$ws.Range("J71").Value = 3500
$ws.Range("L71").Value = 17500
$ws.Range("N71").Value = -24988
$ws.Range("H136").Value = 33912.55
$ws.Range("I136").Value = 5005.9165
$ws.Range("J136").Value = 83466.78999999999
$ws.Range("K136").Value = 15017.7495
$ws.Range("L136").Value = 250400.37
$ws.Range("M136").Value = -12467.7495
$ws.Range("N136").Value = -255500.37

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18186682
$ws.Range("J62").Value = 28575856
$ws.Range("L62").Value = 28575856
$ws.Range("N62").Value = -28577104
$ws.Range("H65").Value = 18186682
$ws.Range("J65").Value = 28575856
$ws.Range("L65").Value = 142879280
$ws.Range("N65").Value = -142885520
$ws.Range("H126").Value = 4483.5
$ws.Range("I126").Value = 3635.625
$ws.Range("K126").Value = 10906.875
$ws.Range("M126").Value = -8436.875
$ws.Range("H132").Value = 8202
$ws.Range("I132").Value = 1284.826
$ws.Range("K132").Value = 3854.478
$ws.Range("M132").Value = -1324.478
$ws.Range("H136").Value = 14909.723
$ws.Range("I136").Value = 1125
$ws.Range("J136").Value = 36571.43
$ws.Range("K136").Value = 3375
$ws.Range("L136").Value = 109714.29
$ws.Range("N136").Value = -114814.29
